# Generate Report for Handback
# The 91f67505-c85e-40b3-925e-b90976650bf8.md file has now been handed back
# (in sync with en-US), so update its status everywhere it is reported, and
# refresh the handback timestamps / clear the stale "old version" error for
# both the zh-cn and de-de target languages.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the 91f67505 file ---
$overview.Range("E3").Value2 = $handedBack
$overview.Range("F3").Value2 = $handedBack

# --- zh-cn sheet: row 3 is the 91f67505 file ---
$zhcn.Range("C3").Value2 = $handedBack
$zhcn.Range("K3").Value2 = "2016-11-08 22:59:25"
$zhcn.Range("P3").Value2 = ""

# --- de-de sheet: row 3 is the 91f67505 file ---
$dede.Range("C3").Value2 = $handedBack
$dede.Range("K3").Value2 = "2016-11-08 22:59:45"
$dede.Range("P3").Value2 = ""

# Error Detail column no longer holds the long stale-version message, so
# shrink it back down to fit its remaining (much shorter) contents.
$zhcn.Columns.Item(16).AutoFit()
$dede.Columns.Item(16).AutoFit()
